# Feat: Added Main screen functions
#
# On the "Translation" sheet:
#  - Normalize the placeholder text in F5/F18/F19 to the plain "<value>"
#    token (previously "Bomb<value>"/"Bomb<value>"/"End<value>").
#  - Append two new rows for the new Main screen functions:
#      SingleUseId19 | Default | Left | LTR | 0
#      SingleUseId20 | Default | Left | LTR | 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("F5").Value = "<value>"
$ws.Range("F18").Value = "<value>"
$ws.Range("F19").Value = "<value>"

$ws.Range("B20").Value = "SingleUseId19"
$ws.Range("C20").Value = "Default"
$ws.Range("D20").Value = "Left"
$ws.Range("E20").Value = "LTR"

# Column F holds text that can look numeric (e.g. "0"); force text storage
# (matching the rest of the column) instead of letting it coerce to a
# number, then drop the number-format style so no stray style id is left
# on the cell.
$f20 = $ws.Range("F20")
$f20.NumberFormat = "@"
$f20.Value = "0"
$f20.Style = "Normal"

$ws.Range("B21").Value = "SingleUseId20"
$ws.Range("C21").Value = "Default"
$ws.Range("D21").Value = "Left"
$ws.Range("E21").Value = "LTR"

$f21 = $ws.Range("F21")
$f21.NumberFormat = "@"
$f21.Value = "0"
$f21.Style = "Normal"
